# Daily attendance processing - 2025-12-18 23:29:15
# Normalize "Recorded By" (column G) values: when the list of recorders
# starts with "System,", rotate it so that entry moves to the end of the
# comma-separated list (swap the first and last items).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G ("Recorded By")
    $val = $cell.Text

    if ($val -ne $null -and $val.StartsWith("System,")) {
        $parts = $val.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        if ($parts.Length -gt 1) {
            $first = $parts[0]
            $last = $parts[$parts.Length - 1]
            $parts[0] = $last
            $parts[$parts.Length - 1] = $first
            $newVal = [string]::Join(", ", $parts)
            $cell.Value = $newVal
        }
    }
}
